$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry row (row 73): Date / Zeit / Einheit / Tätigkeit / Tagebuch note
$ws.Cells.Item(73, 5).Value = 43466
$ws.Cells.Item(73, 6).Value = 3
$ws.Cells.Item(73, 7).Value = "Stunden"
$ws.Cells.Item(73, 8).Value = "Programmieren"
$ws.Cells.Item(73, 9).Value = "Erstellen neuer Klassen zur um die einzelnen Dateien eines Torrents auf der GUI darstellen zu können"
